$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Row 11, column A: was a text "0-02-2015" string, becomes a real date value
# (01-02-2015 => serial 42036), reusing the same date number-format already
# used by the other "Datum" cells (copy format from A10, then set the value).
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A11").Value = (Get-Date -Year 2015 -Month 2 -Day 1 -Hour 0 -Minute 0 -Second 0)

# B11 ("5 uur") and C11 ("multipartconfig in ajax") keep their text as-is.

# New row 12
$ws.Range("A10").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("A12").Value = (Get-Date -Year 2015 -Month 2 -Day 1 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B12").Value = "7 uur"
$ws.Range("C12").Value = "Update organism, errorhandling, dal"

# New row 13
$ws.Range("A10").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = (Get-Date -Year 2015 -Month 3 -Day 2 -Hour 0 -Minute 0 -Second 0)
$ws.Range("C13").Value = "ajax, error, ddl vs jQuery research, subfamily fix, code"
$ws.Range("B13").Value = "10 uur"

$excel.CutCopyMode = 0

# Update view: zoom to 85%, and move the active selection to B13
$ws.Activate()
$excel.ActiveWindow.Zoom = 85
$ws.Range("B13").Select()
